# prodnorm (beef) test with mae
# Fills in the train (MAE) / test (MAE) results (columns X and Y, rows 5-54)
# for the "Random Forest-100 (superdataset-04.csv)" block on the
# "Pilot (normalization)" sheet, and updates the sheet view (zoom/selection)
# to reflect where the author was working when the data was entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pilot (normalization)")

# Data rows: X = train (MAE), Y = test (MAE)
$data = @(
    @(5, 114.2102203679928, 304.72915312232601),
    @(6, 117.6023170731704, 272.11596236099149),
    @(7, 113.2428476679501, 306.7702395209572),
    @(8, 113.77523534445839, 304.90144568006758),
    @(9, 113.5882648694905, 297.98573994867331),
    @(10, 115.82318998716271, 290.73056458511468),
    @(11, 115.3226187419766, 259.32384088964852),
    @(12, 117.9692554557121, 286.34200171086309),
    @(13, 113.0126358579372, 318.42052181351511),
    @(14, 114.0532370560544, 293.32704020530292),
    @(15, 111.51623876765051, 308.41772455089739),
    @(16, 113.6519383825414, 302.01995722839962),
    @(17, 114.272794180573, 299.59416595380588),
    @(18, 111.7445293110823, 311.48528656971712),
    @(19, 116.22065682498901, 295.91679213002482),
    @(20, 112.498305519897, 293.80284003421639),
    @(21, 116.3571031236625, 278.54857142857048),
    @(22, 113.4533119383822, 306.7192899914449),
    @(23, 114.75555840821529, 311.41953806672291),
    @(24, 110.03208600770181, 351.29738237810022),
    @(25, 115.3849850235341, 293.38546621043548),
    @(26, 113.8674475823702, 301.76304533789482),
    @(27, 115.64785836542541, 319.86799828913507),
    @(28, 112.5017501069744, 327.56884516680839),
    @(29, 112.7088232777061, 329.15692899914382),
    @(30, 116.3607766367134, 295.95787852865612),
    @(31, 112.5597496790754, 314.41535500427642),
    @(32, 111.3870068463839, 318.25546621043537),
    @(33, 117.6047347026098, 312.1778956372961),
    @(34, 110.51770860076989, 325.57817792985372),
    @(35, 115.3152952503206, 295.24444824636367),
    @(36, 114.39954642704291, 296.72047048759549),
    @(37, 117.5566816431319, 282.77567151411381),
    @(38, 107.61238125802279, 334.2915055603072),
    @(39, 114.0549657680784, 289.44798973481539),
    @(40, 112.5964591356437, 346.20058169375471),
    @(41, 116.29734060761631, 312.77099230111128),
    @(42, 116.7398973042359, 281.78976903336093),
    @(43, 114.0090479246895, 301.08337040205231),
    @(44, 114.03609114248999, 326.04302822925501),
    @(45, 115.3979396662384, 292.33131736526877),
    @(46, 112.9597304236197, 296.73999999999933),
    @(47, 118.5162644415914, 284.66197604790341),
    @(48, 110.2189473684207, 314.6394781864833),
    @(49, 112.671658108686, 316.1872882805809),
    @(50, 115.0180979888743, 279.00597946963143),
    @(51, 113.0720068463839, 286.11579982891288),
    @(52, 109.6637227214374, 346.06162532078599),
    @(53, 113.7494458707742, 290.71368691188968),
    @(54, 118.34697047496761, 286.63973481608127)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("X$r").Value = $row[1]
    $ws.Range("Y$r").Value = $row[2]
}

# Force a recalculation so the AVERAGE / STDEV.S formulas in X56:Y56 / X57:Y57
# pick up the newly entered values instead of #DIV/0!
$excel.Calculate()

# Update the view: scrolled down and zoomed in further compared to before
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("AA39").Select()
